$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank spacer rows 16:18 (rows that become obsolete, content below shifts up)
$ws.Rows("16:18").Delete()
